# Update for wind and grid
# Wind and Grid items are now green LEDs instead of RGBs. Add a "Pin" column
# with the Arduino pin number used for each item (or "N/A" for the legacy
# "Used" actions), change the Wind/Grid "None"/"Available" actions to the
# same "On"/"Off" actions used by the other items, and leave a comment on
# the old "Used" rows explaining that they're relics from the rgb-led code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Range("D1").Value = "Pin"

# Compressor
$ws.Range("D2").Value = 13
$ws.Range("D3").Value = 13

# Boiler
$ws.Range("D4").Value = 12
$ws.Range("D5").Value = 12

# Reactor
$ws.Range("D6").Value = 8
$ws.Range("D7").Value = 8

# Fischer
$ws.Range("D8").Value = 7
$ws.Range("D9").Value = 7

# Wind - now uses On/Off like the other items, plus a relic "Used" row
$ws.Range("C10").Value = "On"
$ws.Range("D10").Value = 11
$ws.Range("C11").Value = "Off"
$ws.Range("D11").Value = 11
$ws.Range("D12").Value = "N/A"

# Grid - now uses On/Off like the other items, plus a relic "Used" row
$ws.Range("B13").Value = "Grid"
$ws.Range("C13").Value = "On"
$ws.Range("D13").Value = 6
$ws.Range("B14").Value = "Grid"
$ws.Range("C14").Value = "Off"
$ws.Range("D14").Value = 6
$ws.Range("B15").Value = "Grid"
$ws.Range("C15").Value = "Used"
$ws.Range("D15").Value = "N/A"

# Leave explanatory comments on the relic "Used" rows
$commentText = "Author:`nRelic action from when item used rgb led. Can be used in future."
$null = $ws.Range("D12").AddComment($commentText)
$null = $ws.Range("D15").AddComment($commentText)

# Match the author's final selection
$null = $ws.Range("C15").Select()
